$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: India ---
$ws.Range("B14").Value = 96169
$ws.Range("C14").Value = 471
$ws.Range("D14").Value = 36824
$ws.Range("E14").Value = 56316
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 3029

# --- Row 100: Kirguistan ---
$ws.Range("B100").Value = 1216
$ws.Range("C100").Value = 78
$ws.Range("D100").Value = 827
$ws.Range("E100").Value = 375

# --- Rows 154/155: Birmania and Islas Feroe swap order / update stats ---
# Row 154 becomes Birmania (was Islas Feroe)
$ws.Range("A154").Value = "Birmania"
$ws.Range("B154").Value = 187
$ws.Range("C154").Value = 3
$ws.Range("D154").Value = 97
$ws.Range("E154").Value = 84
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 6

# Row 155 becomes Islas Feroe (was Birmania)
$ws.Range("A155").Value = "Islas Feroe"
$ws.Range("B155").Value = 187
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 187
$ws.Range("E155").Value = 0
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 0

# --- Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 06:35"
